$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 5
$ws.Range("J6").Value = "https://www.nytimes.com/2019/12/22/us/arts-health-effects-ucl-study.html"
$ws.Range("B6").Value = "Another Benefit to Going to Museums? You May Live Longer"
$ws.Range("G6").Value = "Mortality"
$ws.Range("H6").Value = "Museum attendance"
$ws.Range("E6").Value = "UK"
$ws.Range("D6").Value = "correlation, causation"
$ws.Range("C6").Value = "culture, health"
$ws.Range("F6").Value = 2019
$ws.Range("L6").Value = 20191225

$ws.Hyperlinks.Add($ws.Range("J6"), "https://www.nytimes.com/2019/12/22/us/arts-health-effects-ucl-study.html")
$ws.Range("J6").Style = "Hyperlink"

$ws.Range("C6").Select()
